# Swap the content of specific columns between row pairs (7,8) and (18,19).
# Columns involved: A, B, E, F, G, H, Q, R, AC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

function Swap-Rows($ws, $cols, $r1, $r2) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

Swap-Rows $ws $cols 7 8
Swap-Rows $ws $cols 18 19
